$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('omnidirectional')
$ws.Range("D3").Value = 37.9
$ws.Range("C4").Value = 39.8
$ws.Range("D4").Value = 42.6

$ws = $wb.Worksheets.Item('345 to 15.0')
$ws.Range("E2").Value = 23.1
$ws.Range("E3").Value = 24.1
$ws.Range("E4").Value = 26.2

$ws = $wb.Worksheets.Item('15.0 to 45.0')
$ws.Range("E2").Value = 21.3
$ws.Range("D3").Value = 25.3
$ws.Range("D4").Value = 29.6

$ws = $wb.Worksheets.Item('45.0 to 75.0')
$ws.Range("D2").Value = 23.8
$ws.Range("E2").Value = 21.6
$ws.Range("C3").Value = 24.1
$ws.Range("D3").Value = 25.7
$ws.Range("E3").Value = 22.9
$ws.Range("D4").Value = 29.9
$ws.Range("E4").Value = 25.5

$ws = $wb.Worksheets.Item('75.0 to 105.0')
$ws.Range("C2").Value = 23.2
$ws.Range("C3").Value = 24.7
$ws.Range("D3").Value = 26
$ws.Range("D4").Value = 30.1

$ws = $wb.Worksheets.Item('105.0 to 135.0')
$ws.Range("D2").Value = 24.6
$ws.Range("C3").Value = 24.9
$ws.Range("D3").Value = 26.3
$ws.Range("C4").Value = 28
$ws.Range("D4").Value = 30.1
$ws.Range("E4").Value = 26.3

$ws = $wb.Worksheets.Item('135.0 to 165.0')
$ws.Range("D2").Value = 25.5
$ws.Range("E2").Value = 23.5
$ws.Range("D3").Value = 27.4
$ws.Range("E3").Value = 24.7
$ws.Range("C4").Value = 29.3
$ws.Range("D4").Value = 31.5
$ws.Range("E4").Value = 27.3

$ws = $wb.Worksheets.Item('165.0 to 195.0')
$ws.Range("D2").Value = 30.8
$ws.Range("E2").Value = 27.9
$ws.Range("D3").Value = 33.5
$ws.Range("E3").Value = 29.5
$ws.Range("D4").Value = 39.4
$ws.Range("E4").Value = 33.1

$ws = $wb.Worksheets.Item('195.0 to 225.0')
$ws.Range("C2").Value = 30.6
$ws.Range("D2").Value = 32
$ws.Range("E2").Value = 29.5
$ws.Range("C3").Value = 32.6
$ws.Range("D3").Value = 34.4
$ws.Range("E3").Value = 31.1
$ws.Range("C4").Value = 36.8
$ws.Range("D4").Value = 39.7
$ws.Range("E4").Value = 34.5

$ws = $wb.Worksheets.Item('225.0 to 255.0')
$ws.Range("D3").Value = 35.6
$ws.Range("E3").Value = 32
$ws.Range("E4").Value = 35.5

$ws = $wb.Worksheets.Item('255.0 to 285.0')
$ws.Range("E2").Value = 30.8
$ws.Range("D3").Value = 36.6
$ws.Range("C4").Value = 39.4
$ws.Range("E4").Value = 36.6

$ws = $wb.Worksheets.Item('285.0 to 315.0')
$ws.Range("E3").Value = 30.3
$ws.Range("C4").Value = 36
$ws.Range("D4").Value = 38.9
$ws.Range("E4").Value = 33.4

$ws = $wb.Worksheets.Item('315.0 to 345.0')
$ws.Range("D2").Value = 29.2
$ws.Range("E2").Value = 27.2
$ws.Range("D3").Value = 31.2
$ws.Range("E3").Value = 28.4
$ws.Range("C4").Value = 32.8
$ws.Range("D4").Value = 35.4
